# add global map method
#
# 1) E14:E23 formulas change from "C-D" to "ABS(C-D)" (values become
#    non-negative). E14 is a standalone formula; E15:E23 is a shared
#    formula group (master E15) - update each independently so the
#    existing shared-formula grouping is preserved.
# 2) A brand new results block ("global map method") is appended at
#    rows 26-36, mirroring the two existing blocks above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new block: "global map method" (rows 26-36) ---------------------------
$ws.Range("A26").Value = "global map method"

$newData = @(
    @(0, 11, 11),
    @(1, 16, 16),
    @(2, 20, 16),
    @(3, 12, 11),
    @(4, 15, 14),
    @(5, 19, 19),
    @(6, 27, 28),
    @(7, 33, 33),
    @(8, 33, 33),
    @(9, 32, 33)
)

$r = 26
foreach ($row in $newData) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

$ws.Cells.Item(26, 5).Formula = "=ABS(C26-D26)"
$ws.Range("E27:E35").Formula = "=ABS(C27-D27)"

# totals row
$ws.Cells.Item(36, 3).Formula = "=SUM(C26:C35)"
$ws.Cells.Item(36, 4).Formula = "=SUM(D26:D35)"
$ws.Cells.Item(36, 5).Formula = "=SUM(E26:E35)"
$ws.Cells.Item(36, 7).Formula = "=C36/D36"

# --- update the "added parameter tuning" block (rows 14-24) ---------------
$ws.Cells.Item(14, 5).Formula = "=ABS(C14-D14)"
$ws.Range("E15:E23").Formula = "=ABS(C15-D15)"

# --- view state: selection moves to G24, view scrolls toward the new block -
[void]$ws.Range("A21").Select()
[void]$ws.Range("G24").Select()
